$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap columns A and B for rows 2-13 (A had year, B had month; now A has month, B has year)
for ($r = 2; $r -le 13; $r++) {
    $a = $ws.Cells.Item($r, 1).Value2
    $b = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 1).Value = $b
    $ws.Cells.Item($r, 2).Value = $a
}

# Add new header columns H1:P1
$ws.Range("H1").Value = "grade_total"
$ws.Range("I1").Value = "grade_distance"
$ws.Range("J1").Value = "grade_visitation"
$ws.Range("K1").Value = "grade_encounters"
$ws.Range("L1").Value = "NEVER"
$ws.Range("M1").Value = "RARELY"
$ws.Range("N1").Value = "SOMETIMES"
$ws.Range("O1").Value = "FREQUENTLY"
$ws.Range("P1").Value = "ALWAYS"

# Add new column data for rows 2-13
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = 0
    $ws.Cells.Item($r, 11).Value = 0
    $ws.Cells.Item($r, 12).Value = 1.017
    $ws.Cells.Item($r, 13).Value = 1.011
    $ws.Cells.Item($r, 14).Value = 1.035
    $ws.Cells.Item($r, 15).Value = 1.121
    $ws.Cells.Item($r, 16).Value = 1.817
}
